# Burndown.xlsx — "Add files via upload" commit re-creation.
#
# The underlying edit adds two sprint-report rows' worth of "Items
# Completed" data (column F) on the Burndown sheet:
#   F11 = "17, 18"
#   F12 = "19, 25"
#   F13 = "26, 29, 30, 31, 32, 33, 34, 35, 36, 37 "   (note trailing space)
#   F14 = 38
# Columns G/H/J/K are formulas (completed-item counts / burndown
# remaining) that recompute automatically once F is populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown")
$ws.Activate()

$ws.Range("F11").Value = "17, 18"
$ws.Range("F12").Value = "19, 25"
$ws.Range("F13").Value = "26, 29, 30, 31, 32, 33, 34, 35, 36, 37 "
$ws.Range("F14").Value = 38

# Matches the author's final cursor position recorded in the workbook.
$ws.Range("F14").Select() | Out-Null
